$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.527.42'
$ws.Range('E2').Value = '  +2.41%  '
$ws.Range('D3').Value = '1.872.75'
$ws.Range('E3').Value = '  +1.65%  '
$ws.Range('E4').Value = '  +0.85%  '
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('E6').Value = '  +0.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4786'
$ws.Range('E7').Value = '  +0.78%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3786'
$ws.Range('E8').Value = '  +3.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07365'
$ws.Range('E9').Value = '  +2.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9393'
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.76'
$ws.Range('E11').Value = '  +5.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07851'
$ws.Range('E12').Value = '  +2.18%  '
$ws.Range('D13').Value = '1.878.51'
$ws.Range('E13').Value = '  +1.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.450'
$ws.Range('E14').Value = '  +2.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.592'
$ws.Range('E15').Value = '  +3.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.98'
$ws.Range('E16').Value = '  +2.74%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.016'
$ws.Range('E17').Value = '  +0.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.97'
$ws.Range('E20').Value = '  +2.95%  '
$ws.Range('D21').Value = '27.557.19'
$ws.Range('E21').Value = '  +2.41%  '
$ws.Range('E22').Value = '  +1.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.75'
$ws.Range('E23').Value = '  +1.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.966'
$ws.Range('E24').Value = '  +2.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.36'
$ws.Range('E25').Value = '  +1.33%  '
$ws.Range('E26').Value = '  +2.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.025'
$ws.Range('E27').Value = '  +1.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '116.01'
$ws.Range('E28').Value = '  +1.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.009'
$ws.Range('E29').Value = '  +1.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08937'
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.336'
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('E32').Value = '  +4.10%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7556'
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.612'
$ws.Range('E34').Value = '  +3.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.709'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02058'
$ws.Range('E36').Value = '  +5.57%  '
$ws.Range('E37').Value = '  +2.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05289'
$ws.Range('E39').Value = '  +1.42%  '
$ws.Range('E40').Value = '  +3.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.093'
$ws.Range('E41').Value = '  +1.95%  '
$ws.Range('E42').Value = '  +1.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.464'
$ws.Range('E43').Value = '  +3.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.64'
$ws.Range('E44').Value = '  +1.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4830'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.662'
$ws.Range('E47').Value = '  +4.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.94'
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '67.42'
$ws.Range('E49').Value = '  +3.23%  '
$ws.Range('E50').Value = '  +1.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9263'
$ws.Range('E51').Value = '  +4.56%  '
